$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The data table runs from row 2 through row 57 (one row per month).
# Append the next month's data as a new row 58, matching the formatting
# (date style) used by the existing rows above it.
$newRow = 58
$prevRow = $newRow - 1

# Copy the prior row's formatting (date number format/style on column A)
# down into the new row before writing the values.
$ws.Range("A" + $prevRow + ":F" + $prevRow).Copy()
$ws.Range("A" + $newRow + ":F" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($newRow, 1).Value = 45536
$ws.Cells.Item($newRow, 2).Value = -0.526
$ws.Cells.Item($newRow, 3).Value = 0.361
$ws.Cells.Item($newRow, 4).Value = -1.506
$ws.Cells.Item($newRow, 5).Value = 0.347
$ws.Cells.Item($newRow, 6).Value = 1.316
